$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "NA" under duplicate_image_filename (column E) for rows 2 through 21
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 5).Value = "NA"
}
